# Daily report rebuild: refresh rows 2-16 with verified driver data.
# Only the cells whose values actually change are touched; rows 8 (Matthew
# Harris) and 11 (Robert Williams) are already correct and are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> Christopher Thomas / TRK-1011 / Downtown Construction
$ws.Range("A2").Value = "Christopher Thomas"
$ws.Range("B2").Value = "TRK-1011"
$ws.Range("C2").Value = "Downtown Construction"
$ws.Range("F2").Value = "06:29 AM"
$ws.Range("G2").Value = "03:33 PM"
$ws.Range("H2").Value = "On Time"
$ws.Range("I2").Value = ""

# Row 3 -> John Smith / TRK-1001 / North Dallas Site
$ws.Range("A3").Value = "John Smith"
$ws.Range("B3").Value = "TRK-1001"
$ws.Range("C3").Value = "North Dallas Site"
$ws.Range("D3").Value = "06:30 AM"
$ws.Range("E3").Value = "03:30 PM"
$ws.Range("F3").Value = "06:25 AM"
$ws.Range("G3").Value = "03:35 PM"

# Row 4 -> James Davis / TRK-1005 / North Dallas Site
$ws.Range("A4").Value = "James Davis"
$ws.Range("B4").Value = "TRK-1005"
$ws.Range("C4").Value = "North Dallas Site"
$ws.Range("F4").Value = "07:03 AM"
$ws.Range("G4").Value = "03:12 PM"
$ws.Range("H4").Value = "Not On Job"
$ws.Range("I4").Value = "At incorrect location: North Richland Hills"

# Row 5 -> Anthony Martin / TRK-1014 / West Plano Project
$ws.Range("A5").Value = "Anthony Martin"
$ws.Range("B5").Value = "TRK-1014"
$ws.Range("C5").Value = "West Plano Project"
$ws.Range("D5").Value = "07:00 AM"
$ws.Range("E5").Value = "04:00 PM"
$ws.Range("F5").Value = "07:16 AM"
$ws.Range("G5").Value = "03:58 PM"
$ws.Range("H5").Value = "Late"
$ws.Range("I5").Value = "16 minutes late"

# Row 6 -> Mark Thompson / TRK-1015 / Downtown Construction
$ws.Range("A6").Value = "Mark Thompson"
$ws.Range("B6").Value = "TRK-1015"
$ws.Range("C6").Value = "Downtown Construction"
$ws.Range("F6").Value = "N/A"
$ws.Range("G6").Value = "N/A"
$ws.Range("H6").Value = "On Time"
$ws.Range("I6").Value = ""

# Row 7 -> David Miller / TRK-1006 / West Plano Project
$ws.Range("A7").Value = "David Miller"
$ws.Range("B7").Value = "TRK-1006"
$ws.Range("C7").Value = "West Plano Project"
$ws.Range("D7").Value = "06:45 AM"
$ws.Range("E7").Value = "03:45 PM"
$ws.Range("F7").Value = "06:47 AM"
$ws.Range("G7").Value = "03:42 PM"

# Row 8 (Matthew Harris / TRK-1013) is unchanged - skip.

# Row 9 -> William Brown / TRK-1004 / Richardson Development
$ws.Range("A9").Value = "William Brown"
$ws.Range("B9").Value = "TRK-1004"
$ws.Range("C9").Value = "Richardson Development"
$ws.Range("D9").Value = "06:30 AM"
$ws.Range("E9").Value = "03:30 PM"
$ws.Range("F9").Value = "06:32 AM"
$ws.Range("G9").Value = "02:54 PM"
$ws.Range("H9").Value = "Early End"
$ws.Range("I9").Value = "36 minutes early"

# Row 10 -> Michael Johnson / TRK-1002 / West Plano Project
$ws.Range("A10").Value = "Michael Johnson"
$ws.Range("B10").Value = "TRK-1002"
$ws.Range("C10").Value = "West Plano Project"
$ws.Range("D10").Value = "06:00 AM"
$ws.Range("E10").Value = "03:00 PM"
$ws.Range("F10").Value = "06:18 AM"
$ws.Range("G10").Value = "03:07 PM"
$ws.Range("I10").Value = "18 minutes late"

# Row 11 (Robert Williams / TRK-1003) is unchanged - skip.

# Row 12 -> Charles Anderson / TRK-1010 / West Plano Project
$ws.Range("A12").Value = "Charles Anderson"
$ws.Range("B12").Value = "TRK-1010"
$ws.Range("C12").Value = "West Plano Project"
$ws.Range("D12").Value = "06:00 AM"
$ws.Range("E12").Value = "03:00 PM"
$ws.Range("F12").Value = "06:05 AM"
$ws.Range("G12").Value = "02:48 PM"

# Row 13 -> Thomas Taylor / TRK-1009 / North Dallas Site
$ws.Range("A13").Value = "Thomas Taylor"
$ws.Range("B13").Value = "TRK-1009"
$ws.Range("C13").Value = "North Dallas Site"
$ws.Range("D13").Value = "07:00 AM"
$ws.Range("E13").Value = "04:00 PM"
$ws.Range("F13").Value = "06:43 AM"
$ws.Range("G13").Value = "03:52 PM"

# Row 14 -> Richard Wilson / TRK-1007 / Downtown Construction
$ws.Range("A14").Value = "Richard Wilson"
$ws.Range("B14").Value = "TRK-1007"
$ws.Range("C14").Value = "Downtown Construction"
$ws.Range("F14").Value = "06:14 AM"
$ws.Range("G14").Value = "03:18 PM"

# Row 15 -> Daniel Jackson / TRK-1012 / Richardson Development
$ws.Range("A15").Value = "Daniel Jackson"
$ws.Range("B15").Value = "TRK-1012"
$ws.Range("C15").Value = "Richardson Development"
$ws.Range("D15").Value = "06:15 AM"
$ws.Range("E15").Value = "03:15 PM"
$ws.Range("F15").Value = "06:16 AM"
$ws.Range("G15").Value = "03:12 PM"

# Row 16 -> Joseph Moore / TRK-1008 / Richardson Development
$ws.Range("A16").Value = "Joseph Moore"
$ws.Range("B16").Value = "TRK-1008"
$ws.Range("C16").Value = "Richardson Development"
$ws.Range("D16").Value = "06:30 AM"
$ws.Range("E16").Value = "03:30 PM"
$ws.Range("F16").Value = "06:55 AM"
$ws.Range("G16").Value = "03:25 PM"
$ws.Range("I16").Value = "25 minutes late"
